$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Valid_Login_Test (sheet1) data ---
# Row 2
$ws1.Range("B2").Value = "APARNA"
$ws1.Range("C2").Value = "MOHANAN"
$ws1.Range("E2").Value = "aparna_123"
# Row 3
$ws1.Range("B3").Value = "NIKHIL"
$ws1.Range("C3").Value = "KANAN"
$ws1.Range("E3").Value = "kanan_90nik"
# Row 4
$ws1.Range("B4").Value = "SELSHIYA"
$ws1.Range("C4").Value = "STEPHEN"
$ws1.Range("E4").Value = "selshiya677"
# Row 5
$ws1.Range("B5").Value = "CALLY"
$ws1.Range("C5").Value = "JOSEPH"
$ws1.Range("E5").Value = "Joe_cal23"

# D2:D5 hold mailto hyperlinks - replace cleanly (delete + recreate) to avoid
# stray relationships, then restore the Hyperlink style.
$ws1.Hyperlinks.Delete()
$ws1.Range("D2").Value = "aparnamohanan@gmail.com"
$ws1.Range("D3").Value = "nikilkanan@gmail.com"
$ws1.Range("D4").Value = "selstephen@hotmail.com"
$ws1.Range("D5").Value = "caljoe@dmail.com"
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:aparnamohanan@gmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("D3"), "mailto:nikilkanan@gmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("D4"), "mailto:selstephen@hotmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("D5"), "mailto:caljoe@dmail.com") | Out-Null
$ws1.Range("D2:D5").Style = "Hyperlink"

# New column F: expected_msg
$ws1.Range("F1").Value = "`${expected_msg}"
$ws1.Range("F2").Value = "reset your password"
$ws1.Range("F3").Value = "reset your password"
$ws1.Range("F4").Value = "reset your password"
$ws1.Range("F5").Value = "reset your password"

# Column widths (subtract the fixed 5/6-character cell-padding offset that the
# host applies when turning ColumnWidth into the stored OOXML width so the
# saved width lands on target)
$ws1.Columns.Item(4).ColumnWidth = 25.88671875 - 0.8333333333333334
$ws1.Columns.Item(6).ColumnWidth = 49.33203125 - 0.8333333333333334
$ws1.Columns.Item(7).ColumnWidth = 15.88671875 - 0.8333333333333334

# Selection matches target
$ws1.Range("E3").Select()

Write-Host "sheet1 done"

# --- Add new sheet: Invalid_Login_Test ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Invalid_Login_Test"

$ws2.Range("A1").Value = "`${test_case}"
$ws2.Range("B1").Value = "`${mailid}"
$ws2.Range("C1").Value = "`${password}"
$ws2.Range("D1").Value = "`${expected_error}"

$ws2.Range("A2").Value = "Invalid_Login_1"
$ws2.Range("B2").Value = "sivabalan@gmail.com"
$ws2.Range("C2").Value = "siva12345"
$ws2.Range("D2").Value = "Incorrect email or password."

$ws2.Range("A3").Value = "Invalid_Login_2"
$ws2.Range("B3").Value = "rasmoh@gmail.com"
$ws2.Range("C3").Value = "ras_890moh"
$ws2.Range("D3").Value = "Incorrect email or password."

$ws2.Range("A4").Value = "Invalid_Login_3"
$ws2.Range("B4").Value = "nirmal@hotmail.com"
$ws2.Range("C4").Value = "surya_456"
$ws2.Range("D4").Value = "Incorrect email or password."

$ws2.Range("A5").Value = "Invalid_Login_4"
$ws2.Range("B5").Value = "lekshmi@hotmail.com"
$ws2.Range("C5").Value = "thangam34"
$ws2.Range("D5").Value = "Incorrect email or password."

$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:sivabalan@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:rasmoh@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "mailto:nirmal@hotmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), "mailto:lekshmi@hotmail.com") | Out-Null
$ws2.Range("B2:B5").Style = "Hyperlink"

# D2 + F2 get the 10pt, vertically centred font used for the message column
$ws2.Range("D2").Font.Size = 10
$ws2.Range("D2").VerticalAlignment = -4108  # xlCenter
$ws2.Range("F2").Font.Size = 10
$ws2.Range("F2").VerticalAlignment = -4108  # xlCenter

$ws2.Columns.Item(1).ColumnWidth = 13.6640625 - 0.8333333333333334
$ws2.Columns.Item(2).ColumnWidth = 19.6640625 - 0.8333333333333334
$ws2.Columns.Item(3).ColumnWidth = 12.5546875 - 0.8333333333333334
$ws2.Columns.Item(4).ColumnWidth = 25.33203125 - 0.8333333333333334
$ws2.Columns.Item(5).ColumnWidth = 14 - 0.8333333333333334
$ws2.Columns.Item(6).ColumnWidth = 28.33203125 - 0.8333333333333334

$ws2.Range("C8").Select()

Write-Host "sheet2 done"

# Restore Valid_Login_Test as the active/selected tab (matches target tabSelected)
$ws1.Activate()
